# 订单导入模板.xlsx -- add "*" (required-field) markers to the order
# import template header row, highlight the merchant-name column, widen
# column A, and move the active selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header label edits (mark required fields with a leading "*") ---
$ws.Range("A1").Value = "*订单编号"
$ws.Range("C1").Value = "*商家ID"
$ws.Range("E1").Value = "*订单履约状态"

# --- Highlight the "商家名称" header (D1) in red ---
$ws.Range("D1").Font.Color = 255

# --- Widen column A ---
# Excel's ColumnWidth property is expressed in "characters of the Normal
# style font" and Excel pads it by 5px (at 7px/char for the default font)
# when it serialises the column's stored width; 100/7 characters round-
# trips to a stored width of exactly 15.
$ws.Columns("A").ColumnWidth = 100 / 7

# --- Move the active selection ---
[void]$ws.Range("G19").Select()
